$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.448.03'
$ws.Range("E2").Value = '  +0.66%  '

$ws.Range("D3").Value = '2.249.75'
$ws.Range("E3").Value = '  -0.05%  '

$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.08'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.34%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.50'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.58%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.572'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.71%  '

$ws.Range("E8").Value = '  +0.21%  '

$ws.Range("E9").Value = '  -0.87%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.86'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.50%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0813'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.81%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.24'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.70%  '

$ws.Range("E13").Value = '  +0.69%  '

$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.841'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.75%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '2.242.62'
$ws.Range("E15").Value = '  -0.52%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.72'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.28%  '

$ws.Range("D17").Value = '44.122.27'
$ws.Range("E17").Value = '  +0.29%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '0.0₃0965'
$ws.Range("E18").Value = '  -0.75%  '

$ws.Range("B19").Value = 'InternetComputer(DFINITY)'
$ws.Range("C19").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.47'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.41'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '65.98'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '237.87'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.52%  '

$ws.Range("E23").Value = '  +2.36%  '

$ws.Range("E24").Value = '  +2.63%  '

$ws.Range("E25").Value = '  -0.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '38.52'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +5.08%  '

$ws.Range("E27").Value = '  +4.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.90'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.85%  '

$ws.Range("E29").Value = '  -3.12%  '

$ws.Range("E30").Value = '  +0.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '154.49'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.37%  '

$ws.Range("E32").Value = '  -1.83%  '

$ws.Range("E33").Value = '  -0.39%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.12'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -10.67%  '

$ws.Range("E35").Value = '  +2.58%  '

$ws.Range("E36").Value = '  +0.69%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.82'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.49'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.57%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.81'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.25%  '

$ws.Range("E40").Value = '  -0.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0305'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.35%  '

$ws.Range("E42").Value = '  +0.28%  '

$ws.Range("D43").Value = '1.751.50'
$ws.Range("E43").Value = '  -0.51%  '

$ws.Range("E44").Value = '  +1.32%  '

$ws.Range("E45").Value = '  -6.68%  '

$ws.Range("B46").Value = 'THORChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.97'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.52%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '99.95'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.22%  '

$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.62'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +6.79%  '

$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '70.90'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.51%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '56.14'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.44%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.17'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.84%  '
